$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'26.391.54"
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Formula = "'1.833.23"
$ws.Range("E3").Value = '  -0.11%  '
$ws.Range("D4").Formula = "'1.001"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Formula = "'254.63"
$ws.Range("E5").Value = '  -1.84%  '
$ws.Range("D7").Formula = "'0.5287"
$ws.Range("E7").Value = '  +0.54%  '
$ws.Range("D8").Formula = "'0.2830"
$ws.Range("E8").Value = '  -11.31%  '
$ws.Range("D9").Formula = "'0.06899"
$ws.Range("E9").Value = '  +1.82%  '
$ws.Range("D10").Formula = "'1.847.83"
$ws.Range("E10").Value = '  +0.11%  '
$ws.Range("D11").Formula = "'16.57"
$ws.Range("E11").Value = '  -11.34%  '
$ws.Range("D12").Formula = "'0.6981"
$ws.Range("E12").Value = '  -10.41%  '
$ws.Range("D13").Formula = "'0.07135"
$ws.Range("E13").Value = '  -7.66%  '
$ws.Range("D14").Formula = "'87.00"
$ws.Range("E14").Value = '  -0.56%  '
$ws.Range("E15").Value = '  -2.40%  '
$ws.Range("D16").Formula = "'1.001"
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("D18").Formula = "'13.26"
$ws.Range("E18").Value = '  -3.88%  '
$ws.Range("D19").Formula = "'0.000007386"
$ws.Range("E19").Value = '  -6.78%  '
$ws.Range("D20").Formula = "'26.418.44"
$ws.Range("E20").Value = '  -0.34%  '
$ws.Range("D21").Formula = "'2.085.76"
$ws.Range("E21").Value = '  +0.96%  '
$ws.Range("D22").Formula = "'4.509"
$ws.Range("E22").Value = '  -1.96%  '
$ws.Range("D23").Formula = "'5.824"
$ws.Range("E23").Value = '  -2.25%  '
$ws.Range("E24").Value = '  -3.85%  '
$ws.Range("D25").Formula = "'142.22"
$ws.Range("E25").Value = '  +0.47%  '
$ws.Range("D26").Formula = "'1.680"
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("E27").Value = '  -5.07%  '
$ws.Range("D28").Formula = "'16.63"
$ws.Range("E28").Value = '  -1.68%  '
$ws.Range("D29").Formula = "'109.12"
$ws.Range("E29").Value = '  -2.06%  '
$ws.Range("D30").Formula = "'4.098"
$ws.Range("E30").Value = '  -0.76%  '
$ws.Range("D31").Formula = "'0.08735"
$ws.Range("E31").Value = '  +0.69%  '
$ws.Range("D32").Formula = "'3.877"
$ws.Range("E32").Value = '  -4.32%  '
$ws.Range("D33").Formula = "'0.04707"
$ws.Range("E33").Value = '  -2.70%  '
$ws.Range("D34").Formula = "'2.886"
$ws.Range("E34").Value = '  +1.22%  '
$ws.Range("D35").Formula = "'1.112"
$ws.Range("E35").Value = '  -1.32%  '
$ws.Range("D36").Formula = "'0.7070"
$ws.Range("E36").Value = '  -2.83%  '
$ws.Range("D37").Formula = "'3.061"
$ws.Range("E37").Value = '  -0.82%  '
$ws.Range("D38").Formula = "'2.189"
$ws.Range("E38").Value = '  -1.64%  '
$ws.Range("D39").Formula = "'0.01645"
$ws.Range("E39").Value = '  -6.69%  '
$ws.Range("D40").Formula = "'0.4485"
$ws.Range("E40").Value = '  -5.22%  '
$ws.Range("D41").Formula = "'0.8649"
$ws.Range("E41").Value = '  -2.96%  '
$ws.Range("D42").Formula = "'105.15"
$ws.Range("E42").Value = '  -3.77%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").Formula = "'5.740"
$ws.Range("E44").Value = '  -2.65%  '
$ws.Range("D45").Formula = "'7.058"
$ws.Range("E45").Value = '  -7.21%  '
$ws.Range("D46").Formula = "'8.672"
$ws.Range("E46").Value = '  -3.27%  '
$ws.Range("D47").Formula = "'0.1190"
$ws.Range("E47").Value = '  -2.67%  '
$ws.Range("D48").Formula = "'33.48"
$ws.Range("E48").Value = '  -3.49%  '
$ws.Range("D49").Formula = "'58.79"
$ws.Range("E49").Value = '  -1.12%  '
$ws.Range("D50").Formula = "'0.05576"
$ws.Range("E50").Value = '  -4.70%  '
$ws.Range("D51").Formula = "'0.8615"
$ws.Range("E51").Value = '  -3.37%  '
